$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates per diff; numeric-looking text values are forced to Text
# format first so Excel stores the literal string instead of
# auto-converting it to a floating point number.

$ws.Range('D2').Value = '35.483.21'
$ws.Range('E2').Value = '  +0.72%  '
$ws.Range('D3').Value = '1.924.51'
$ws.Range('E3').Value = '  +1.71%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.727'
$ws.Range('E5').Value = '  +11.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '254.15'
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '40.95'
$ws.Range('E8').Value = '  -0.66%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.356'
$ws.Range('E9').Value = '  +2.65%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '52.52'
$ws.Range('E10').Value = '  +5.08%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0746'
$ws.Range('E11').Value = '  +5.46%  '
$ws.Range('E12').Value = '  +0.66%  '
$ws.Range('D13').Value = '2.203.32'
$ws.Range('E13').Value = '  +1.68%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '12.78'
$ws.Range('E14').Value = '  +7.49%  '
$ws.Range('E15').Value = '  +4.18%  '
$ws.Range('D16').Value = '1.934.36'
$ws.Range('E16').Value = '  +2.11%  '
$ws.Range('E17').Value = '  +1.86%  '
$ws.Range('D18').Value = '35.475.54'
$ws.Range('E18').Value = '  +0.75%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '74.48'
$ws.Range('E19').Value = '  +4.78%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '243.58'
$ws.Range('E21').Value = '  +1.25%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '13.03'
$ws.Range('E22').Value = '  +5.37%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.10'
$ws.Range('E23').Value = '  +8.06%  '
$ws.Range('E24').Value = '  -0.03%  '
$ws.Range('E25').Value = '  +2.46%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.40'
$ws.Range('E26').Value = '  -0.72%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '167.82'
$ws.Range('E27').Value = '  -1.34%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.64'
$ws.Range('E28').Value = '  +3.09%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.134'
$ws.Range('E29').Value = '  +6.63%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '18.75'
$ws.Range('E30').Value = '  +3.01%  '
$ws.Range('D31').Value = '4.125.99'
$ws.Range('E31').Value = '  +19.38%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.39'
$ws.Range('E32').Value = '  +7.13%  '
$ws.Range('B33').Value = 'TrustWalletToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.65'
$ws.Range('E33').Value = '  +25.25%  '
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.98'
$ws.Range('E34').Value = '  +14.44%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0582'
$ws.Range('E35').Value = '  +4.27%  '
$ws.Range('E36').Value = '  +3.72%  '
$ws.Range('E37').Value = '  +0.12%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.919'
$ws.Range('E38').Value = '  -1.86%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.03'
$ws.Range('E39').Value = '  +0.94%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.46'
$ws.Range('E40').Value = '  +9.82%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '97.30'
$ws.Range('E41').Value = '  +9.52%  '
$ws.Range('E42').Value = '  +4.31%  '
$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0658'
$ws.Range('E43').Value = '  +3.56%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0211'
$ws.Range('E44').Value = '  +1.52%  '
$ws.Range('D45').Value = '1.347.79'
$ws.Range('E45').Value = '  +0.87%  '
$ws.Range('E46').Value = '  +4.78%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.81'
$ws.Range('E47').Value = '  +4.56%  '
$ws.Range('E48').Value = '  +0.86%  '
$ws.Range('E49').Value = '  +0.58%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '45.30'
$ws.Range('E50').Value = '  -5.29%  '
$ws.Range('E51').Value = '  +6.13%  '
